$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'58.275.72"
$ws.Range("E2").Value = "  -2.65%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.574.77"
$ws.Range("E3").Value = "  -2.85%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'536.31"
$ws.Range("E5").Value = "  -0.11%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'143.15"
$ws.Range("E6").Value = "  -1.54%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.576"
$ws.Range("E8").Value = "  +0.75%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "'6.79"
$ws.Range("E9").Value = "  +1.48%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0998"
$ws.Range("E10").Value = "  -3.64%  "

# Row 11 - was Cardano, now TRON
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.138"
$ws.Range("E11").Value = "  +2.80%  "

# Row 12 - was TRON, now Cardano
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  -2.19%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'3.039.29"
$ws.Range("E13").Value = "  -2.81%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "'58.187.87"
$ws.Range("E14").Value = "  -2.67%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'20.58"
$ws.Range("E15").Value = "  -3.03%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'2.579.79"
$ws.Range("E16").Value = "  -3.38%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -2.05%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "'4.46"
$ws.Range("E18").Value = "  +0.89%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'334.14"
$ws.Range("E19").Value = "  -3.04%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'10.01"
$ws.Range("E20").Value = "  -2.68%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.11"
$ws.Range("E21").Value = "  -3.87%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.09%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'66.56"
$ws.Range("E23").Value = "  -0.68%  "

# Row 24 - Polygon
$ws.Range("D24").Value = "'0.417"
$ws.Range("E24").Value = "  +0.55%  "

# Row 25 - Binance-PegBSC-USD
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - Kaspa
$ws.Range("D26").Value = "'0.158"
$ws.Range("E26").Value = "  -4.89%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'7.03"
$ws.Range("E27").Value = "  -4.29%  "

# Row 28 - was PEPE, now USDe
$ws.Range("B28").Value = "USDe"
$ws.Range("C28").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.07%  "

# Row 29 - was USDe, now PEPE
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0₃0731"
$ws.Range("E29").Value = "  -2.74%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'1.63"
$ws.Range("E30").Value = "  -1.65%  "

# Row 31 - Aptos
$ws.Range("D31").Value = "'5.93"
$ws.Range("E31").Value = "  +1.11%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'155.14"
$ws.Range("E32").Value = "  +3.25%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'18.78"
$ws.Range("E33").Value = "  -1.58%  "

# Row 34 - NEARProtocol
$ws.Range("D34").Value = "'3.88"
$ws.Range("E34").Value = "  -4.24%  "

# Row 35 - OKB
$ws.Range("D35").Value = "'36.94"
$ws.Range("E35").Value = "  -0.56%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'1.09"
$ws.Range("E36").Value = "  -5.60%  "

# Row 37 - SuiNetwork
$ws.Range("D37").Value = "'0.842"
$ws.Range("E37").Value = "  +2.28%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "'0.815"
$ws.Range("E38").Value = "  -3.53%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -4.26%  "

# Row 40 - Filecoin
$ws.Range("D40").Value = "'3.57"
$ws.Range("E40").Value = "  -0.74%  "

# Row 41 - Bittensor
$ws.Range("D41").Value = "'277.83"
$ws.Range("E41").Value = "  -5.71%  "

# Row 42 - FirstDigitalUSD
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.19%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -2.36%  "

# Row 44 - WhiteBITCoin
$ws.Range("E44").Value = "  -0.89%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "'0.0531"
$ws.Range("E45").Value = "  -2.81%  "

# Row 46 - Stellar
$ws.Range("D46").Value = "'0.0937"
$ws.Range("E46").Value = "  -1.97%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "'18.44"
$ws.Range("E47").Value = "  -4.75%  "

# Row 48 - VeChain
$ws.Range("D48").Value = "'0.0225"
$ws.Range("E48").Value = "  -0.84%  "

# Row 49 - Maker
$ws.Range("D49").Value = "'1.909.84"
$ws.Range("E49").Value = "  -3.32%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'17.83"
$ws.Range("E50").Value = "  -3.53%  "

# Row 51 - RenderToken
$ws.Range("D51").Value = "'4.38"
$ws.Range("E51").Value = "  -4.11%  "

# The leading apostrophes above force text entry (so numeric-looking
# strings like "1.00" or date-like "3.039.29" aren't coerced into
# numbers/dates), but that also stamps a "quote prefix" cell style.
# Reset the style on the whole price column back to Normal so no stray
# formatting is introduced versus the original (unstyled) cells.
$ws.Range("D2:D51").Style = "Normal"
